$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 133 (the "よろこべ" post), shifting all subsequent rows up by one.
$ws.Rows.Item(133).Delete()
